# Fruta / hortaliza, semanal
#
# A new weekly price-survey record is inserted as row 374 on the
# "Hortaliza, Terminal La Palmera de La Serena - Papa" sheet. Every
# existing record from row 374 onward shifts down by one row (374->375,
# ..., 463->464), and the new row 374 is populated with a fresh date
# (D374) and fresh volume/price figures (J/K/L/M/P), while inheriting
# the rest of the record's fields (variety, quality, unit, origin, etc.)
# from the row that is being pushed down beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row above row 374; everything below shifts down
#    by one (old 374 -> 375, ..., old 463 -> 464).
$ws.Rows.Item(374).Insert()

# 2) Seed the new row 374 by duplicating the record that now sits in
#    row 375 (the original row 374), so all of its non-numeric fields
#    (variety, quality, unit, origin, classification, etc.) carry over.
$ws.Rows.Item(375).Copy()
$ws.Rows.Item(374).PasteSpecial()

# 3) Overwrite the new record's date and volume/price columns with the
#    actual reported values for this entry.
$ws.Cells.Item(374, 4).Value = 44785    # D374 - Fecha
$ws.Cells.Item(374, 10).Value = 2000    # J374 - Volumen
$ws.Cells.Item(374, 11).Value = 12500   # K374 - Precio minimo
$ws.Cells.Item(374, 12).Value = 13000   # L374 - Precio maximo
$ws.Cells.Item(374, 13).Value = 12750   # M374 - Precio promedio ponderado
$ws.Cells.Item(374, 16).Value = 510     # P374 - Precio $/Kg
